# PROS-10194 - CCRU - Promo tracking KPIs
#
# Renames/expands a handful of shared "scene tag" strings used throughout the
# Canteen KPI sheet:
#   - "Panoramic Photo"                                    -> "Panoramic Photo, SS_Panoramic Photo"
#   - "Panoramic photo of Cooler"                          -> "Panoramic photo of Cooler, SS_Panoramic photo of Cooler - Horeca"   (only the "Scenes to exclude" cell, X38)
#   - "Menu Board, Cash Zone, SS_Menu Board, SS_Cash Zone" -> "Menu Board, Cash Zone, SS_Cash Zone - Canteen, QSR, SS_Menu Board - Canteen, QSR"
#
# All other cells keep their original text; Excel's shared-string table is
# renumbered automatically as the old "Panoramic Photo" entry disappears and
# the new/edited strings are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Canteen")

# --- "Scenes to include" (column Y) cells that used the bare "Panoramic Photo" tag ---
$panoramicPhotoRows = @(4,5,6,7,8,9,10,11,12,13,14,16,17,18,20,21,22,23,25,26,27,28,29,30,31)
foreach ($r in $panoramicPhotoRows) {
    $cell = $ws.Cells.Item($r, 25)  # column Y
    if ($cell.Value2 -eq "Panoramic Photo") {
        $cell.Value2 = "Panoramic Photo, SS_Panoramic Photo"
    }
}

# --- "Scenes to exclude" (column X), row 38: Cooler prime-position KPI ---
$cellX38 = $ws.Cells.Item(38, 24)  # column X
if ($cellX38.Value2 -eq "Panoramic photo of Cooler") {
    $cellX38.Value2 = "Panoramic photo of Cooler, SS_Panoramic photo of Cooler - Horeca"
}

# --- "Scenes to exclude" (column X) cells for the Activation KPIs ---
$menuBoardRows = @(43,44,45,46,47,48,49)
foreach ($r in $menuBoardRows) {
    $cell = $ws.Cells.Item($r, 24)  # column X
    if ($cell.Value2 -eq "Menu Board, Cash Zone, SS_Menu Board, SS_Cash Zone") {
        $cell.Value2 = "Menu Board, Cash Zone, SS_Cash Zone - Canteen, QSR, SS_Menu Board - Canteen, QSR"
    }
}

# --- Cursor / selection moved to AD5 before saving ---
$ws.Range("AD5").Select()
